# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1160840
$ws.Range("C4").Value = 66
$ws.Range("E4").Value = 919667

# Row 16 - Belgica
$ws.Range("B16").Value = 49906
$ws.Range("C16").Value = 389
$ws.Range("D16").Value = 12309
$ws.Range("E16").Value = 29753
$ws.Range("F16").Value = 674
$ws.Range("G16").Value = 79
$ws.Range("H16").Value = 7844

# Row 39 - Indonesia
$ws.Range("B39").Value = 11192
$ws.Range("C39").Value = 349
$ws.Range("D39").Value = 1876
$ws.Range("E39").Value = 8471
$ws.Range("G39").Value = 14
$ws.Range("H39").Value = 845

# Row 78 - Estonia
$ws.Range("B78").Value = 1700
$ws.Range("C78").Value = 1
$ws.Range("D78").Value = 259
$ws.Range("E78").Value = 1386
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 55

# Row 85 - Eslovenia
$ws.Range("D85").Value = 241
$ws.Range("E85").Value = 1102
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 96

# Row 91 - Hong Kong
$ws.Range("D91").Value = 879
$ws.Range("E91").Value = 157
$ws.Range("F91").Value = 3
